$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: entry for 1/23/2020 ---
$ws.Range("A13").Value = 43853
$ws.Range("A13").NumberFormat = "m/d/yy"

$ws.Range("B13").Value = "class"
$ws.Range("C13").Value = "Nic"
$ws.Range("D13").Value = "Search for certain functions, and write down thought process and priority."
$ws.Range("E13").Value = "Was able to locate and identify potential findings and solutioins to our problem. "
$ws.Range("F13").Value = "The activity allowed us to be very purposeful and reflective for our every move. It forced us to really slow down and think about what decisions we were making, and make a usually thoughtless process to become very meaningful."
$ws.Range("G13").Value = "To be very honest this activity felt….very slow, and I questioned whether or not it would be useful in the future. The first minute into the activity felt meaningful, but it quickly felt like a manual task that took more time trying to find reasons to my actions than getting the task at hand done."

# --- Row 14: entry for 1/25/2020 ---
$ws.Range("A14").Value = 43855
$ws.Range("A14").NumberFormat = "mm-dd-yy"

$ws.Range("B14").Value = "3:00 - 7:00"
$ws.Range("C14").Value = "Nic, Rafael, Chris"
$ws.Range("D14").Value = "Finding 2 features to look for in FreeCol and creating a UML document."
$ws.Range("E14").Value = "We were able to successfully identify 2 features we wanted to look for, find them within the program, and create a uml documentation of the entire project."
$ws.Range("F14").Value = "There was so much code! Luckily the code was not as hard to sift through, thanks to the project being well documented. But after having created the UML document, the code base seemed so vast and big that it made me think that perhaps we had almost gotten lucky in terms of how relatively quickly we got through the initial process."
$ws.Range("G14").Value = "Even though we had a few struggles, I realized very quickly how in comparison to some other programs, FreeCol is probably quite small in size. This made me appreciate people's ability to actually ""read"" code. Even with nice documentation, it may still be very hard to find or understand certain programs...so imagine a program that has poor to no documentation! just thinking about it gives me chills."

# --- Update the view to reflect the last edited cell ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("G14").Select()
